$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PctPopUrbanGDPPerCapRealUSD")

$newColO = @(
    -2.5638043529732282, -1.9333293013201871, -2.5015872292678432, -2.595009255364269, -2.7904622454973569, -2.2128553639228921, -1.727101771729634, 1.1896713589043659,
    -2.0599921269573329, -0.58885984623404464, -2.4849063860753509, -1.528034942999013, -0.58521144759678034, -2.4146984047149802, -1.0794431210827149, -0.79047591821232599,
    -0.69620251389204446, -1.4387138410440981, -2.0814111831035502, -2.17876494731996, -2.222137521538802, -1.498301796825479, -2.0602729846840679, -5.2703563903582751,
    -1.951162019417475, -3.324317226211376, -2.1063300609021458, -2.071912764624495, -1.1084819198036899, -1.308767918691963, -1.8192268771293709, -3.1675980229825562,
    -1.378615170960996, -2.1197420544188161, -0.59645923001078116, -1.968824181252067, -2.0573093779077372, -0.99079914278383074, -7.2453164054770776, -2.166306004480218,
    -1.923017801874819, -2.1908010983117299, -2.748641818731179, -2.4506051342417101, -18.040670518023258, -2.1213339977736392, -1.557857847788215, -1.708218596908089,
    -0.81815043665245191, -3.4765063281912112, -2.000421963789123, -0.83575044073620497, -1.882344554530639, -0.57679142768196268, -3.173196083845732, -0.86795587885461012,
    -2.416895584734418, -1.552935003874216, -0.31960062222834362, -2.3039521631149542, -1.884701478940177, -2.8000208103686219, -0.46892191352523233, -1.5292953896156209,
    -2.4766821776990469, -1.4166974892347779, -3.276025622554219, -2.150949543969837, 0.2279088621277465, -1.8835200148229341, -2.5527275254083501, -3.1139172834254851,
    -2.0347814564705229, -1.344125740527178, -3.5531123571160199, -1.258331271851314, -3.3039302312256331, -1.9961549048662059, -2.3457625400034261, -5.6232983257116684,
    -2.3046135189108452, -1.178253605740319, -2.2471789866673308, -1.922956133221897, -1.4982244064000521, -1.219490598262444, -2.5101645918006712, -3.514929380728069,
    -1.761995001543794, -2.2991768927130232, -1.9340209619469131, -1.896026509636713, -1.5834650294452299, -2.1518223022991601, -1.0385124687186651, -1.4551041270718681,
    -1.881905889801972, -3.006878062117893, -16.128387224167788, -2.234557968804844, -1.3789843028809601, -1.602659131324522, -2.807380829680707, -1.5142493606996721,
    -3.1557623644588721, -2.131620143027297, -1.4024509000222489, -1.1073908193552751, -1.836818521304274, -1.3843378691300201, -1.3285220151799779, -1.152911611315341,
    -2.6341214830080371, -2.112233498473584, -2.2776806176397701, -2.619289920597569, -3.9450832163289009, -2.0358653776459068, -1.339914478594993, -6.0613666365157206,
    -2.1994665822257979, -2.056949524355498, 0.31924638912835213, -2.1948338957660471, -1.496797327314574, -1.4580332051604461, -1.1303461262420711, -5.8972549820172491,
    -2.1514519230217051, -2.1690494492313031, -12.124795192440899, -1.881898909100997, -2.7585320079900182, -2.4996639170257571, -0.78220285780793741, -4.8164369576977446,
    -3.560036637817773, -1.5924099781167029, -2.5160300971999372, -2.1436228820040348, -1.888889008994302, -3.538723149636597, -0.34066288367558573, -2.2091532010369912,
    -1.5685355309425619, -1.756281665942053
)

for ($i = 0; $i -lt $newColO.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 15).Value = $newColO[$i]
}
